$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for rows 2-19 (columns A=Name, B=Position, C=Team)
$data = @(
    @(2,  "Tyus Jones",            "PG",      "Phoenix Suns"),
    @(3,  "Anfernee Simons",       "PG,SG",   "Portland Trail Blazers"),
    @(4,  "Fred VanVleet",         "PG",      "Houston Rockets"),
    @(5,  "James Harden",          "PG,SG",   "LA Clippers"),
    @(6,  "Anthony Edwards",       "SG,SF",   "Minnesota Timberwolves"),
    @(7,  "Paul George",           "SG,SF,PF","Philadelphia 76ers"),
    @(8,  "Jayson Tatum",          "SF,PF",   "Boston Celtics"),
    @(9,  "Andrew Wiggins",        "SF,PF",   "Golden State Warriors"),
    @(10, "Jaren Jackson Jr.",     "PF,C",    "Memphis Grizzlies"),
    @(11, "Bobby Portis",          "PF,C",    "Milwaukee Bucks"),
    @(12, "Jaden Ivey",            "PG,SG",   "Detroit Pistons"),
    @(13, "Ivica Zubac",           "C",       "LA Clippers"),
    @(14, "Aaron Gordon",          "PF,C",    "Denver Nuggets"),
    @(15, "Dillon Brooks",         "SG,SF",   "Houston Rockets"),
    @(16, "Dennis Schröder",       "PG",      "Golden State Warriors"),
    @(17, "Zion Williamson",       "PF,C",    "New Orleans Pelicans"),
    @(18, "Kyle Kuzma",            "PF",      "Washington Wizards"),
    @(19, "Giannis Antetokounmpo", "PF,C",    "Milwaukee Bucks")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
